$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) holds the same "changed" date serial number (45192 -> 45202)
# for every data row from row 2 through row 459. Update them all in one shot.
$ws.Range("C2:C459").Value = 45202
